$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the obsolete "BTemp" EEPROM-data-sections row (row 40: C=10, BTemp, 32). ---
# This shifts every row below it up by one (old r41->40, r42->41, r43->42) and the
# shared-string table drops the now-unused "BTemp" entry, re-indexing later strings.
$ws.Rows(40).Delete()

# --- Re-establish the running-total formulas in column C for the EEPROM section. ---
$ws.Range("C32").Formula = "=C31+E32"
$ws.Range("C33").Formula = "=C32+E33"
$ws.Range("C34").Formula = "=C33+E34"
$ws.Range("C35").Formula = "=C34+E35"
$ws.Range("C36").Formula = "=C35+E36"
$ws.Range("C37").Formula = "=C36+E37"
$ws.Range("C38").Formula = "=C37+E38"
$ws.Range("C39").Formula = "=C38+E39"

# Row 40 (now "BTempSet") restarts the byte offset at 32 rather than continuing the chain.
$ws.Range("C40").Formula = "=32"

# Row 41 (now "BHistSet") is a fixed literal, no longer a formula.
$ws.Range("C41").Value = 74

# Row 42 (now "HTempSet") resumes the running total from row 41.
$ws.Range("C42").Formula = "=C41+E42"

# --- Center-align the running-total column for the EEPROM data section rows. ---
$ws.Range("C33:C42").HorizontalAlignment = -4108

$ws.Range("C42").Select()
